# Auto update Excel log
# Appends new sensor-log rows captured at 2026-01-28 16:39 to the PIR,
# Humidity and Temperature sheets.
#
# Column A (dates like "2026-01-28") and the Humidity sheet's Value column
# (percentages like "87.0%") look numeric/date-like to Excel, which would
# otherwise silently convert them to a date serial / fraction on write. A
# leading apostrophe forces those specific values to stay plain text, which
# matches how every pre-existing row in these logs is stored. Columns that
# Excel already keeps as text on their own (times, "Bathroom", "Active",
# "Inactive", "No Motion", "22.8C", ...) are written as-is.

function Write-LogRows {
    # NOTE: called positionally everywhere below - binding array arguments
    # through a *named* parameter (-Rows $x) unwraps the outer array in this
    # PowerShell host, so $Rows.Count comes back 0. Positional binding keeps
    # the array-of-arrays intact.
    param($Worksheet, $StartRow, $Rows, $TextForceCols)
    for ($i = 0; $i -lt $Rows.Count; $i++) {
        $r = $StartRow + $i
        $row = $Rows[$i]
        for ($c = 0; $c -lt $row.Count; $c++) {
            $col = $c + 1
            if ($TextForceCols -contains $col) {
                $Worksheet.Cells.Item($r, $col).Value = "'" + $row[$c]
            } else {
                $Worksheet.Cells.Item($r, $col).Value = $row[$c]
            }
        }
    }
}

$wb = $excel.ActiveWorkbook

# ---- PIR sheet: add rows 47-60 -------------------------------------------------
$ws = $wb.Worksheets.Item("PIR")
$pirRows = @(
    @("2026-01-28","16:39:00","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","16:39:00","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","16:39:02","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","16:39:07","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","16:39:12","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","16:39:17","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","16:39:22","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","16:39:27","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","16:39:32","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","16:39:37","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","16:39:42","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","16:39:48","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","16:39:52","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","16:39:57","16:00","Bathroom","No Motion","Inactive")
)
Write-LogRows $ws 47 $pirRows @(1)

# ---- Humidity sheet: add rows 48-59 --------------------------------------------
$ws = $wb.Worksheets.Item("Humidity")
$humidityRows = @(
    @("2026-01-28","16:39:00","16:00","Bathroom","87.0%","Active"),
    @("2026-01-28","16:39:01","16:00","Bathroom","87.9%","Active"),
    @("2026-01-28","16:39:07","16:00","Bathroom","87.0%","Active"),
    @("2026-01-28","16:39:11","16:00","Bathroom","88.0%","Active"),
    @("2026-01-28","16:39:23","16:00","Bathroom","87.9%","Active"),
    @("2026-01-28","16:39:27","16:00","Bathroom","87.9%","Active"),
    @("2026-01-28","16:39:35","16:00","Bathroom","87.0%","Active"),
    @("2026-01-28","16:39:39","16:00","Bathroom","87.9%","Active"),
    @("2026-01-28","16:39:47","16:00","Bathroom","87.0%","Active"),
    @("2026-01-28","16:39:51","16:00","Bathroom","87.9%","Active"),
    @("2026-01-28","16:39:55","16:00","Bathroom","87.0%","Active"),
    @("2026-01-28","16:39:59","16:00","Bathroom","87.9%","Active")
)
Write-LogRows $ws 48 $humidityRows @(1,5)

# ---- Temperature sheet: add rows 48-59 -----------------------------------------
$ws = $wb.Worksheets.Item("Temperature")
$temperatureRows = @(
    @("2026-01-28","16:39:00","16:00","Bathroom","22.8C","Active"),
    @("2026-01-28","16:39:01","16:00","Bathroom","22.8C","Active"),
    @("2026-01-28","16:39:07","16:00","Bathroom","22.8C","Active"),
    @("2026-01-28","16:39:11","16:00","Bathroom","22.9C","Active"),
    @("2026-01-28","16:39:23","16:00","Bathroom","22.8C","Active"),
    @("2026-01-28","16:39:27","16:00","Bathroom","22.8C","Active"),
    @("2026-01-28","16:39:35","16:00","Bathroom","22.8C","Active"),
    @("2026-01-28","16:39:39","16:00","Bathroom","22.8C","Active"),
    @("2026-01-28","16:39:47","16:00","Bathroom","22.8C","Active"),
    @("2026-01-28","16:39:51","16:00","Bathroom","22.8C","Active"),
    @("2026-01-28","16:39:55","16:00","Bathroom","22.8C","Active"),
    @("2026-01-28","16:39:59","16:00","Bathroom","22.8C","Active")
)
Write-LogRows $ws 48 $temperatureRows @(1)
